$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in F23: nu=8585 -> nu=0.8585
$ws.Range("F23").Value = "nu=0.8585, n_components=3, freq bands (Hz) 4-8,8-13,13-30"

# Add new row 27
$ws.Range("A27").Value = "Bandpower + PCA + NuSVM (linear kernel)"
$ws.Range("B27").Value = 0.8433
$ws.Range("B27").NumberFormat = "0.00%"
$ws.Range("C27").Value = "17/19"
$ws.Range("D27").Value = "L"
$ws.Range("E27").Value = "0, 1, 1, 2, 3, 3, 5, 12, 13, 23, 30, 52, 57"
$ws.Range("F27").Value = "nu=0.8585, n_components=3, freq bands (Hz) 4-8,8-13,13-30"

$ws.Range("D27").Select()
